# Apply the "Test cases and Pom pages added" edit:
#  - Insert a new "CountryAndState" worksheet right after "Login" (becomes sheetId 8,
#    pushing Place/Register/SearchBox/Subscribe's r:id's up by one).
#  - Populate it with Country/State headers and a United States/Texas data row.
#  - Make CountryAndState the active/selected sheet (tabSelected), with C9 selected.
#  - Login sheet keeps its own new selection (C2), loses tabSelected, and both
#    sheets pick up best-fit column widths.

$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("Login")

# Login: change the remembered selection to C2 (still just a plain selection,
# not the active tab anymore once CountryAndState is selected below).
$login.Range("C2").Select()

# Give Login's used columns a best-fit width (A:B).
$login.Range("A1:B2").EntireColumn.AutoFit()

# Insert the new sheet right after "Login" so it lands in the 3rd tab position.
$newSheet = $wb.Worksheets.Add($null, $login)
$newSheet.Name = "CountryAndState"

$newSheet.Range("A1").Value = "Country"
$newSheet.Range("B1").Value = "State"
$newSheet.Range("A2").Value = "United States"
$newSheet.Range("B2").Value = "Texas"

# Best-fit column A on the new sheet too (only column A got a <col> entry
# in the authored file - column B was left at the default width).
$newSheet.Range("A1:A2").EntireColumn.AutoFit()

# Leave the new sheet as the active tab/selected cell, matching the diff.
$newSheet.Range("C9").Select()
